# This script applies the following changes described by the diff:
#  1) Sheet "展览" (exhibitions): small updates to a handful of "want to go" counts (column F).
#  2) Sheet "演出" (performances): small updates to a handful of "want to go" counts (column F).
#  3) Sheet "本地生活" (local life): small updates to a handful of "want to go" counts (column F).
#  4) Sheet "全部类型" (all types, an aggregate of the above three sheets): the duplicated
#     row 7 (a repeat of the "夏川里美" concert entry) is removed, which shifts every row below
#     it up by one and shrinks the used range from A1:I43 to A1:I42. After the shift, the
#     "want to go" counts (column F) for the rows whose events were also refreshed in the
#     source sheets above are updated to match the newly refreshed values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet "展览"
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F9").Value2  = 1463
$wsExpo.Range("F11").Value2 = 1347
$wsExpo.Range("F12").Value2 = 3006
$wsExpo.Range("F13").Value2 = 427
$wsExpo.Range("F14").Value2 = 1626
$wsExpo.Range("F16").Value2 = 797
$wsExpo.Range("F18").Value2 = 1390
$wsExpo.Range("F20").Value2 = 65
$wsExpo.Range("F21").Value2 = 1125
$wsExpo.Range("F23").Value2 = 3
$wsExpo.Range("F24").Value2 = 3488
$wsExpo.Range("F25").Value2 = 688
$wsExpo.Range("F27").Value2 = 1544

# ---------------------------------------------------------------------------
# 2) Sheet "演出"
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F6").Value2  = 10
$wsShow.Range("F7").Value2  = 52
$wsShow.Range("F8").Value2  = 21
$wsShow.Range("F12").Value2 = 81
$wsShow.Range("F13").Value2 = 16

# ---------------------------------------------------------------------------
# 3) Sheet "本地生活"
# ---------------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value2 = 796
$wsLocal.Range("F3").Value2 = 6

# ---------------------------------------------------------------------------
# 4) Sheet "全部类型"
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

# Remove the duplicated row (row 7), which shifts all following rows up by one
# and shrinks the sheet's dimension from A1:I43 to A1:I42.
$wsAll.Rows.Item(7).Delete()

# After the shift, refresh the "want to go" counts (column F) for the rows whose
# underlying events were also refreshed in the "展览"/"演出"/"本地生活" sheets above.
$wsAll.Range("F3").Value2  = 796
$wsAll.Range("F5").Value2  = 6
$wsAll.Range("F10").Value2 = 10
$wsAll.Range("F11").Value2 = 52
$wsAll.Range("F13").Value2 = 21
$wsAll.Range("F19").Value2 = 1463
$wsAll.Range("F21").Value2 = 1347
$wsAll.Range("F22").Value2 = 3006
$wsAll.Range("F23").Value2 = 427
$wsAll.Range("F24").Value2 = 1626
$wsAll.Range("F26").Value2 = 797
$wsAll.Range("F28").Value2 = 1390
$wsAll.Range("F30").Value2 = 65
$wsAll.Range("F33").Value2 = 1125
$wsAll.Range("F35").Value2 = 3
$wsAll.Range("F36").Value2 = 3488
$wsAll.Range("F37").Value2 = 688
$wsAll.Range("F39").Value2 = 1544
$wsAll.Range("F40").Value2 = 81
$wsAll.Range("F41").Value2 = 16
